# Updates odds values in the Betfair "Jogos do Dia" workbook to match the
# newer snapshot of the data (commit: "Atualizando o arquivo XLSX").
# Only numeric odd cell values change; no structural changes are made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.95
$ws.Range("H2").Value = 2.74
$ws.Range("I3").Value = 1.68
$ws.Range("Q3").Value = 1.34
$ws.Range("K5").Value = 3.95
$ws.Range("F6").Value = 1.76
$ws.Range("K6").Value = 4.4
$ws.Range("J8").Value = 4.5
$ws.Range("F9").Value = 1.98
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1.88
$ws.Range("F10").Value = 2.84
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 2.54
$ws.Range("I10").Value = 2.68
$ws.Range("J10").Value = 3.6
$ws.Range("F12").Value = 1.83
$ws.Range("G12").Value = 1.86
$ws.Range("H12").Value = 4.6
$ws.Range("J12").Value = 3.9
$ws.Range("P12").Value = 2.26
$ws.Range("Q12").Value = 1.72
$ws.Range("F13").Value = 1.3
$ws.Range("G13").Value = 1.31
$ws.Range("H13").Value = 13.5
$ws.Range("I13").Value = 15.5
$ws.Range("K13").Value = 6.4
$ws.Range("P13").Value = 2.3
$ws.Range("S13").Value = 2.78
$ws.Range("T13").Value = 2.26
$ws.Range("X13").Value = 26
$ws.Range("Y13").Value = 80
$ws.Range("AB13").Value = 8.800000000000001
$ws.Range("AD13").Value = 160
$ws.Range("AF13").Value = 7.4
$ws.Range("AH13").Value = 75
$ws.Range("AJ13").Value = 9.6
$ws.Range("AK13").Value = 15.5
$ws.Range("AL13").Value = 130
$ws.Range("AN13").Value = 5.1
$ws.Range("G14").Value = 1.42
$ws.Range("H14").Value = 9.4
$ws.Range("J15").Value = 3.85
$ws.Range("K15").Value = 4.1
$ws.Range("Q15").Value = 1.63
$ws.Range("Q16").Value = 1.78
$ws.Range("G17").Value = 1.52
$ws.Range("K17").Value = 5.3
$ws.Range("O17").Value = 1.2
$ws.Range("P17").Value = 2.52
$ws.Range("R17").Value = 1.6
$ws.Range("Z17").Value = 65
$ws.Range("AG17").Value = 10.5
$ws.Range("AI17").Value = 85
$ws.Range("AK17").Value = 16
$ws.Range("AM17").Value = 95
$ws.Range("AN17").Value = 5.9
$ws.Range("G18").Value = 1.96
$ws.Range("I18").Value = 4.4
$ws.Range("K19").Value = 4.6
$ws.Range("P19").Value = 2.06
$ws.Range("Q19").Value = 1.87
$ws.Range("S19").Value = 3.15
$ws.Range("T19").Value = 1.92
$ws.Range("U19").Value = 2
$ws.Range("X19").Value = 17.5
$ws.Range("Z19").Value = 110
$ws.Range("AC19").Value = 9.6
$ws.Range("AF19").Value = 9.6
$ws.Range("AH19").Value = 24
$ws.Range("AJ19").Value = 16.5
$ws.Range("G20").Value = 2.42
$ws.Range("H20").Value = 3.15
$ws.Range("P21").Value = 2.44
$ws.Range("X21").Value = 26
$ws.Range("AJ21").Value = 11
$ws.Range("F22").Value = 1.67
$ws.Range("G22").Value = 1.69
$ws.Range("H22").Value = 5.5
$ws.Range("I22").Value = 5.8
$ws.Range("P22").Value = 2.48
$ws.Range("Q22").Value = 1.58
$ws.Range("F23").Value = 1.78
$ws.Range("Q23").Value = 1.78
$ws.Range("G24").Value = 2.04
$ws.Range("H24").Value = 3.7
$ws.Range("J24").Value = 4
$ws.Range("P25").Value = 2.7
$ws.Range("G26").Value = 2.9
$ws.Range("F29").Value = 2.06
$ws.Range("Q30").Value = 2.88
$ws.Range("F31").Value = 1.86
$ws.Range("P31").Value = 1.91
$ws.Range("T31").Value = 1.81
$ws.Range("F32").Value = 1.85
$ws.Range("G32").Value = 2.14
$ws.Range("H32").Value = 3.9
$ws.Range("I32").Value = 6.4
$ws.Range("J32").Value = 3
$ws.Range("K32").Value = 3.9
